$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H69").Value = 1546.2593
$ws.Range("J69").Value = 1517.2693
$ws.Range("L69").Value = 4551.8079
$ws.Range("N69").Value = -6299.8079
$ws.Range("H72").Value = 1546.2593
$ws.Range("J72").Value = 1517.2693
$ws.Range("L72").Value = 13655.4237
$ws.Range("N72").Value = -22391.4237
$ws.Range("H76").Value = 2927059.5
$ws.Range("I76").Value = 3183.9285
$ws.Range("J76").Value = 11113911
$ws.Range("K76").Value = 3183.9285
$ws.Range("L76").Value = 11113911
$ws.Range("M76").Value = -2868.9285
$ws.Range("N76").Value = -11114541
$ws.Range("H79").Value = 2927059.5
$ws.Range("I79").Value = 3183.9285
$ws.Range("J79").Value = 11113911
$ws.Range("K79").Value = 3183.9285
$ws.Range("L79").Value = 11113911
$ws.Range("M79").Value = -2091.9285
$ws.Range("N79").Value = -11116095
$ws.Range("H132").Value = 3766.6365
$ws.Range("I132").Value = 4003
$ws.Range("J132").Value = 2703
$ws.Range("K132").Value = 12009
$ws.Range("L132").Value = 8109
$ws.Range("M132").Value = -9479
$ws.Range("N132").Value = -13169

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H37").Value = 2000
$ws.Range("J37").Value = 0
$ws.Range("L37").Value = 0
$ws.Range("H80").Value = 48631.25
$ws.Range("J80").Value = 48631.25
$ws.Range("L80").Value = 48631.25
$ws.Range("N80").Value = -50627.25
$ws.Range("H83").Value = 48631.25
$ws.Range("J83").Value = 48631.25
$ws.Range("L83").Value = 145893.75
$ws.Range("N83").Value = -155877.75
$ws.Range("H92").Value = 19000
$ws.Range("J92").Value = 19000
$ws.Range("L92").Value = 19000
$ws.Range("N92").Value = -23992
$ws.Range("H97").Value = 1035.2354
$ws.Range("I97").Value = 1035.2354
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 1035.2354
$ws.Range("L97").Value = 0
$ws.Range("M97").Value = -539.2354
$ws.Range("N37").ClearContents()
$ws.Range("N97").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1318.9565
$ws.Range("I99").Value = 1135.5
$ws.Range("J99").Value = 1519.091
$ws.Range("K99").Value = 1135.5
$ws.Range("L99").Value = 1519.091
$ws.Range("M99").Value = 362.5
$ws.Range("N99").Value = -4515.091

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 400.16666
$ws.Range("I22").Value = 425.25
$ws.Range("J22").Value = 350
$ws.Range("K22").Value = 425.25
$ws.Range("L22").Value = 350
$ws.Range("M22").Value = -75.25
$ws.Range("N22").Value = -1050
$ws.Range("H58").Value = 32599.375
$ws.Range("I58").Value = 1439.3334
$ws.Range("J58").Value = 500000
$ws.Range("K58").Value = 1439.3334
$ws.Range("L58").Value = 500000
$ws.Range("M58").Value = -1236.3334
$ws.Range("N58").Value = -500406
$ws.Range("H136").Value = 32599.375
$ws.Range("I136").Value = 1439.3334
$ws.Range("J136").Value = 500000
$ws.Range("K136").Value = 4318.0002
$ws.Range("L136").Value = 1500000
$ws.Range("M136").Value = -1768.0002
$ws.Range("N136").Value = -1505100

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1727.909
$ws.Range("I5").Value = 1444.6666
$ws.Range("K5").Value = 4333.9998
$ws.Range("M5").Value = -4221.9998
$ws.Range("H34").Value = 853.5714
$ws.Range("I34").Value = 701
$ws.Range("J34").Value = 914.6
$ws.Range("K34").Value = 2103
$ws.Range("L34").Value = 2743.8
$ws.Range("M34").Value = -2019
$ws.Range("N34").Value = -2911.8
$ws.Range("H39").Value = 2269.3333
$ws.Range("I39").Value = 1900
$ws.Range("J39").Value = 2454
$ws.Range("K39").Value = 5700
$ws.Range("L39").Value = 7362
$ws.Range("M39").Value = -5406
$ws.Range("N39").Value = -7950
$ws.Range("H55").Value = 3000
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 3000
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 9000
$ws.Range("N55").Value = -9354
$ws.Range("H56").Value = 3087.1428
$ws.Range("I56").Value = 3087.1428
$ws.Range("K56").Value = 3087.1428
$ws.Range("M56").Value = -2557.1428
$ws.Range("H122").Value = 430.55554
$ws.Range("I122").Value = 234.625
$ws.Range("K122").Value = 2111.625
$ws.Range("M122").Value = 338.375
$ws.Range("H131").Value = 781.39
$ws.Range("I131").Value = 595
$ws.Range("J131").Value = 787.15466
$ws.Range("K131").Value = 1785
$ws.Range("L131").Value = 2361.46398
$ws.Range("M131").Value = 3255
$ws.Range("N131").Value = -12441.46398
$ws.Range("H135").Value = 1727.909
$ws.Range("I135").Value = 1444.6666
$ws.Range("K135").Value = 13001.9994
$ws.Range("M135").Value = -10466.9994
$ws.Range("M55").ClearContents()

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 6157.143
$ws.Range("I22").Value = 7000
$ws.Range("J22").Value = 5525
$ws.Range("K22").Value = 7000
$ws.Range("L22").Value = 5525
$ws.Range("M22").Value = -6705
$ws.Range("N22").Value = -6115
$ws.Range("H27").Value = 6157.143
$ws.Range("I27").Value = 7000
$ws.Range("J27").Value = 5525
$ws.Range("K27").Value = 7000
$ws.Range("L27").Value = 5525
$ws.Range("M27").Value = -6893
$ws.Range("N27").Value = -5739
$ws.Range("H93").Value = 2898.7144
$ws.Range("I93").Value = 2898.7144
$ws.Range("K93").Value = 2898.7144
$ws.Range("M93").Value = -1650.7144
$ws.Range("H136").Value = 42776.918
$ws.Range("I136").Value = 42776.918
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 128330.754
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -125780.754
$ws.Range("N136").ClearContents()

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 14000
$ws.Range("I54").Value = 5000
$ws.Range("J54").Value = 15800
$ws.Range("K54").Value = 5000
$ws.Range("L54").Value = 15800
$ws.Range("M54").Value = -4480
$ws.Range("N54").Value = -16840
$ws.Range("H81").Value = 83334510
$ws.Range("I81").Value = 1284
$ws.Range("K81").Value = 2568
$ws.Range("M81").Value = -1507
$ws.Range("H84").Value = 83334510
$ws.Range("I84").Value = 1284
$ws.Range("K84").Value = 12840
$ws.Range("M84").Value = -7536
$ws.Range("H96").Value = 5550
$ws.Range("I96").Value = 1500
$ws.Range("J96").Value = 9600
$ws.Range("K96").Value = 1500
$ws.Range("L96").Value = 9600
$ws.Range("M96").Value = -127
$ws.Range("N96").Value = -12346
$ws.Range("H107").Value = 1977422.6
$ws.Range("I107").Value = 888.2727
$ws.Range("K107").Value = 2664.8181
$ws.Range("M107").Value = -744.8181
$ws.Range("H132").Value = 2525
$ws.Range("I132").Value = 1426
$ws.Range("J132").Value = 3074.5
$ws.Range("K132").Value = 4278
$ws.Range("L132").Value = 9223.5
$ws.Range("M132").Value = -1748
$ws.Range("N132").Value = -14283.5
$ws.Range("H136").Value = 28573040
$ws.Range("I136").Value = 40001452
$ws.Range("J136").Value = 2009.8
$ws.Range("K136").Value = 120004356
$ws.Range("L136").Value = 6029.4
$ws.Range("M136").Value = -120001806
$ws.Range("N136").Value = -11129.4
